$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.458.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.802.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.602'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.77%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.287'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0666'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.065.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.804.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.630'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.446.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '238.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0766'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.71%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.121'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.61%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0512'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.639'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.300.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0185'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.01%  '
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("E41").Value = '  +1.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.945'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.83%  '
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.79%  '
$ws.Range("E46").Value = '  +2.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.965.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.22%  '
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0612'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.79%  '
